$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4736548282512274
$ws.Range("C2").Value = 0.2167045972234583
$ws.Range("D2").Value = 0.03549638714461167
$ws.Range("F2").Value = 0.6267297712044311
$ws.Range("G2").Value = 0.4675799249437631
$ws.Range("H2").Value = 0.6281267015680569
$ws.Range("I2").Value = 0.6475669555290722
$ws.Range("K2").Value = 0.3035872364264662
$ws.Range("L2").Value = 0.3015734276458062
$ws.Range("O2").Value = 2.149556015944952
$ws.Range("B3").Value = 0.4261442327014038
$ws.Range("C3").Value = 0.2178611245974764
$ws.Range("D3").Value = 0.03238533503401442
$ws.Range("F3").Value = 0.6288625145977704
$ws.Range("G3").Value = 0.4711170883961273
$ws.Range("H3").Value = 0.633257725356394
$ws.Range("I3").Value = 0.654474542775894
$ws.Range("K3").Value = 0.2649351231132471
$ws.Range("L3").Value = 0.2901931964065341
$ws.Range("O3").Value = 2.167564603810419
$ws.Range("B4").Value = 0.3969783963001134
$ws.Range("C4").Value = 0.2186236695307926
$ws.Range("D4").Value = 0.03046043976647184
$ws.Range("F4").Value = 0.6305680382662757
$ws.Range("G4").Value = 0.4736201968416864
$ws.Range("H4").Value = 0.6366780062700741
$ws.Range("I4").Value = 0.6590698188082804
$ws.Range("K4").Value = 0.2411067174988517
$ws.Range("L4").Value = 0.2833747787254453
$ws.Range("O4").Value = 2.179881040412099
$ws.Range("B5").Value = 0.3850953945984088
$ws.Range("C5").Value = 0.2189476328174322
$ws.Range("D5").Value = 0.02967237630085862
$ws.Range("F5").Value = 0.6313626603166611
$ws.Range("G5").Value = 0.4747234837178382
$ws.Range("H5").Value = 0.6381396941951323
$ws.Range("I5").Value = 0.6610314406230131
$ws.Range("K5").Value = 0.2313730080356748
$ws.Range("L5").Value = 0.2806388413277716
$ws.Range("O5").Value = 2.18521668437802
$ws.Range("B6").Value = 0.3831223956892131
$ws.Range("C6").Value = 0.2190022261793061
$ws.Range("D6").Value = 0.02954129958363438
$ws.Range("F6").Value = 0.6315006235317711
$ws.Range("G6").Value = 0.4749117102144922
$ws.Range("H6").Value = 0.6383865085456861
$ws.Range("I6").Value = 0.6613625422681402
$ws.Range("K6").Value = 0.2297553367560567
$ws.Range("L6").Value = 0.2801871184021536
$ws.Range("O6").Value = 2.186121785241113
$ws.Range("B7").Value = 0.3968181275021152
$ws.Range("C7").Value = 0.2186279850389248
$ws.Range("D7").Value = 0.03044982640270888
$ws.Range("F7").Value = 0.6305783514750303
$ws.Range("G7").Value = 0.4736347391484301
$ws.Range("H7").Value = 0.6366974441316415
$ws.Range("I7").Value = 0.6590959135352037
$ws.Range("K7").Value = 0.2409755393205302
$ws.Range("L7").Value = 0.2833377081989994
$ws.Range("O7").Value = 2.179951716904526
$ws.Range("B8").Value = 0.4572724467923024
$ws.Range("C8").Value = 0.2170925067429117
$ws.Range("D8").Value = 0.03442677135549843
$ws.Range("F8").Value = 0.6273829643241342
$ws.Range("G8").Value = 0.4687307509222265
$ws.Range("H8").Value = 0.629839908889366
$ws.Range("I8").Value = 0.6498752189581012
$ws.Range("K8").Value = 0.2902802668666311
$ws.Range("L8").Value = 0.2976144875670457
$ws.Range("O8").Value = 2.155504044476743
$ws.Range("B9").Value = 0.5758386103598241
$ws.Range("C9").Value = 0.2144958996641328
$ws.Range("D9").Value = 0.04210740397899571
$ws.Range("F9").Value = 0.6242584414869725
$ws.Range("G9").Value = 0.4617452119428975
$ws.Range("H9").Value = 0.6185309833406549
$ws.Range("I9").Value = 0.6346023920368005
$ws.Range("K9").Value = 0.3861815770474379
$ws.Range("L9").Value = 0.3269500514772687
$ws.Range("O9").Value = 2.117553894841322
$ws.Range("B10").Value = 0.6629253498201706
$ws.Range("C10").Value = 0.2128387197776789
$ws.Range("D10").Value = 0.0476767791613355
$ws.Range("F10").Value = 0.6238783049092476
$ws.Range("G10").Value = 0.4582210439599166
$ws.Range("H10").Value = 0.6115232417342824
$ws.Range("I10").Value = 0.6250940906750664
$ws.Range("K10").Value = 0.4561352670904739
$ws.Range("L10").Value = 0.349317902736999
$ws.Range("O10").Value = 2.095766085073663
$ws.Range("B11").Value = 0.7025316770110237
$ws.Range("C11").Value = 0.2121387859658839
$ws.Range("D11").Value = 0.05019415081019929
$ws.Range("F11").Value = 0.6241213970008275
$ws.Range("G11").Value = 0.4569678903981895
$ws.Range("H11").Value = 0.6086171248282852
$ws.Range("I11").Value = 0.621140464518259
$ws.Range("K11").Value = 0.4878441844321912
$ws.Range("L11").Value = 0.3596704841940408
$ws.Range("O11").Value = 2.08717825935824
$ws.Range("B12").Value = 0.7175274055467185
$ws.Range("C12").Value = 0.2118814590594482
$ws.Range("D12").Value = 0.05114505207154707
$ws.Range("F12").Value = 0.6242732620389546
$ws.Range("G12").Value = 0.4565437499295157
$ws.Range("H12").Value = 0.6075571226137129
$ws.Range("I12").Value = 0.6196967965971467
$ws.Range("K12").Value = 0.4998346420470909
$ws.Range("L12").Value = 0.3636161673138929
$ws.Range("O12").Value = 2.084116648891197
$ws.Range("B13").Value = 0.7142979250912163
$ws.Range("C13").Value = 0.2119365360377969
$ws.Range("D13").Value = 0.0509403647678397
$ws.Range("F13").Value = 0.6242378952586449
$ws.Range("G13").Value = 0.4566328534277133
$ws.Range("H13").Value = 0.6077836131641234
$ws.Range("I13").Value = 0.6200053371413148
$ws.Range("K13").Value = 0.4972530509119792
$ws.Range("L13").Value = 0.3627652666903742
$ws.Range("O13").Value = 2.084767551470222
$ws.Range("B14").Value = 0.7037654364800972
$ws.Range("C14").Value = 0.2121174609294982
$ws.Range("D14").Value = 0.05027242985906355
$ws.Range("F14").Value = 0.6241326923914272
$ws.Range("G14").Value = 0.4569319855817042
$ws.Range("H14").Value = 0.6085291066917407
$ws.Range("I14").Value = 0.6210206210235931
$ws.Range("K14").Value = 0.4888309918381424
$ws.Range("L14").Value = 0.3599945898257317
$ws.Range("O14").Value = 2.086922562178785
$ws.Range("B15").Value = 0.6973136529645672
$ws.Range("C15").Value = 0.2122292874581859
$ws.Range("D15").Value = 0.04986299000739791
$ws.Range("F15").Value = 0.624076041434094
$ws.Range("G15").Value = 0.4571217785218238
$ws.Range("H15").Value = 0.6089910138798658
$ws.Range("I15").Value = 0.6216494777912764
$ws.Range("K15").Value = 0.4836699995048548
$ws.Range("L15").Value = 0.3583007727074516
$ws.Range("O15").Value = 2.088267368358132
$ws.Range("B16").Value = 0.6603367068794057
$ws.Range("C16").Value = 0.2128855445108186
$ws.Range("D16").Value = 0.04751193359936678
$ws.Range("F16").Value = 0.6238707886592749
$ws.Range("G16").Value = 0.4583099896318785
$ws.Range("H16").Value = 0.6117188298412373
$ws.Range("I16").Value = 0.6253599523279689
$ws.Range("K16").Value = 0.4540606720047435
$ws.Range("L16").Value = 0.3486448957839343
$ws.Range("O16").Value = 2.096353960517249
$ws.Range("B17").Value = 0.637649373402752
$ws.Range("C17").Value = 0.2133019259981168
$ws.Range("D17").Value = 0.04606546207809004
$ws.Range("F17").Value = 0.6238514171396261
$ws.Range("G17").Value = 0.4591286148948441
$ws.Range("H17").Value = 0.6134643936608342
$ws.Range("I17").Value = 0.6277314415719921
$ws.Range("K17").Value = 0.4358667601603372
$ws.Range("L17").Value = 0.3427666659217579
$ws.Range("O17").Value = 2.101653880277993
$ws.Range("B18").Value = 0.6245993446916884
$ws.Range("C18").Value = 0.2135464952191484
$ws.Range("D18").Value = 0.04523197127058154
$ws.Range("F18").Value = 0.623879433522724
$ws.Range("G18").Value = 0.4596324067493853
$ws.Range("H18").Value = 0.6144949164408757
$ws.Range("I18").Value = 0.6291304484974489
$ws.Range("K18").Value = 0.4253914754667676
$ws.Range("L18").Value = 0.3394023620362674
$ws.Range("O18").Value = 2.104826828466258
$ws.Range("B19").Value = 0.6201807067286325
$ws.Range("C19").Value = 0.2136301752068093
$ws.Range("D19").Value = 0.04494950624514615
$ws.Range("F19").Value = 0.6238956452076678
$ws.Range("G19").Value = 0.4598086373597994
$ws.Range("H19").Value = 0.6148483896904366
$ws.Range("I19").Value = 0.6296101367582168
$ws.Range("K19").Value = 0.421842920674095
$ws.Range("L19").Value = 0.3382661380113205
$ws.Range("O19").Value = 2.105922527322647
$ws.Range("B20").Value = 0.6400645758858445
$ws.Range("C20").Value = 0.2132570761960437
$ws.Range("D20").Value = 0.04621959896415717
$ws.Range("F20").Value = 0.6238494266064194
$ws.Range("G20").Value = 0.4590380611284957
$ws.Range("H20").Value = 0.6132758308617525
$ws.Range("I20").Value = 0.627475371058253
$ws.Range("K20").Value = 0.4378046382909133
$ws.Range("L20").Value = 0.3433906861292826
$ws.Range("O20").Value = 2.101076801140977
$ws.Range("B21").Value = 0.7068591541124363
$ws.Range("C21").Value = 0.2120641095566356
$ws.Range("D21").Value = 0.05046868331331211
$ws.Range("F21").Value = 0.6241619698598342
$ws.Range("G21").Value = 0.4568427547731488
$ws.Range("H21").Value = 0.6083090387839576
$ws.Range("I21").Value = 0.6207209557643374
$ws.Range("K21").Value = 0.4913052224122794
$ws.Range("L21").Value = 0.3608077167601635
$ws.Range("O21").Value = 2.086284414868132
$ws.Range("B22").Value = 0.7504992823737098
$ws.Range("C22").Value = 0.2113294371017389
$ws.Range("D22").Value = 0.05323185296151678
$ws.Range("F22").Value = 0.6247148549045889
$ws.Range("G22").Value = 0.4557018062885803
$ws.Range("H22").Value = 0.6052988932662373
$ws.Range("I22").Value = 0.6166183160617322
$ws.Range("K22").Value = 0.52617140215321
$ws.Range("L22").Value = 0.3723386631575352
$ws.Range("O22").Value = 2.077726645206099
$ws.Range("B23").Value = 0.727209304553071
$ws.Range("C23").Value = 0.2117174384403029
$ws.Range("D23").Value = 0.05175838077848027
$ws.Range("F23").Value = 0.6243878751601954
$ws.Range("G23").Value = 0.4562838451065758
$ws.Range("H23").Value = 0.6068838868094133
$ws.Range("I23").Value = 0.6187794368862392
$ws.Range("K23").Value = 0.5075720112177748
$ws.Range("L23").Value = 0.3661708846576914
$ws.Range("O23").Value = 2.082192503037106
$ws.Range("B24").Value = 0.6389726838885394
$ws.Range("C24").Value = 0.2132773366352012
$ws.Range("D24").Value = 0.04614991957717507
$ws.Range("F24").Value = 0.6238502045667857
$ws.Range("G24").Value = 0.4590788972221489
$ws.Range("H24").Value = 0.6133609960755209
$ws.Range("I24").Value = 0.6275910296291727
$ws.Range("K24").Value = 0.4369285713861473
$ws.Range("L24").Value = 0.343108519337548
$ws.Range("O24").Value = 2.101337306289338
$ws.Range("B25").Value = 0.5437651131650227
$ws.Range("C25").Value = 0.2151542002648128
$ws.Range("D25").Value = 0.04004239885826877
$ws.Range("F25").Value = 0.6247673471357373
$ws.Range("G25").Value = 0.4633528758287682
$ws.Range("H25").Value = 0.6213616842069101
$ws.Range("I25").Value = 0.6384333549791137
$ws.Range("K25").Value = 0.3603245194033491
$ws.Range("L25").Value = 0.3188707898409007
$ws.Range("O25").Value = 2.126750425267261
